$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target state for rows 2..19 (columns A..E) after:
#  - inserting a new leading data row (date 39400 / 2007 -> 2008), which
#    pushes all the previously-existing data rows down by one, and
#  - recomputing the "forecast" columns (C = y_0_forecast, E = y_1_forecast)
#    for every row, per the bugfix described in the commit message.
$data = @(
    @(39400, 2007, 1.144978573787081,   2008, $null),
    @(39765, 2008, 1.381024225294869,   2009, $null),
    @(40130, 2009, -0.8792832172735965, 2010, $null),
    @(40494, 2010, 0.9337833426867448,  2011, 2.776721259569048),
    @(40862, 2011, 2.791140000794279,   2012, 1.397663935193227),
    @(41228, 2012, 0.4451370000809973,  2013, 0.1571217018783377),
    @(41592, 2013, 0.2545814083968478,  2014, 1.564494423159446),
    @(41957, 2014, 1.297015177357297,   2015, 0.4654670911270031),
    @(42321, 2015, 1.365576377841027,   2016, 1.13561893405616),
    @(42689, 2016, 2.204449574611278,   2017, 1.120909053655028),
    @(43053, 2017, 2.18621550610123,    2018, 2.012342079251583),
    @(43418, 2018, 1.911050033324102,   2019, 1.897633937626764),
    @(43783, 2019, 1.457852003181337,   2020, 0.5018989042238298),
    @(44159, 2020, -3.258619210312885,  2021, -1.154010671959604),
    @(44525, 2021, 0.4255262881966981,  2022, 2.711235034366544),
    @(44890, 2022, 3.293290997728171,   2023, 0.8446287801135322),
    @(45254, 2023, -0.2814561130375925, 2024, 0.2005659760944134),
    @(45618, 2024, -0.6470065423293869, 2025, 1.323978403813442)
)

$firstDataRow = 2
$lastOldDataRow = 19  # previous last row (18) shifted down by one

# Shift the bordered/centered date values in column A down by one row first
# (bottom-up, so a value is never overwritten before it has been read).
for ($r = $lastOldDataRow - 1; $r -ge $firstDataRow; $r--) {
    $ws.Cells.Item($r + 1, 1).Value = $ws.Cells.Item($r, 1).Value
}

# Write the final values for every data row/column.
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + $firstDataRow
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]

    if ($row[2] -eq $null) { $ws.Cells.Item($r, 3).ClearContents() }
    else { $ws.Cells.Item($r, 3).Value = $row[2] }

    $ws.Cells.Item($r, 4).Value = $row[3]

    if ($row[4] -eq $null) { $ws.Cells.Item($r, 5).ClearContents() }
    else { $ws.Cells.Item($r, 5).Value = $row[4] }
}

# Row 2 (brand-new row) and row 19 (newly extended last row) need the same
# bordered / centered / date-formatted style that the rest of column A
# already has. Copy that formatting (not the value) from an existing,
# correctly-styled date cell.
$ws.Cells.Item(3, 1).Copy() | Out-Null
$ws.Cells.Item($firstDataRow, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item($lastOldDataRow, 1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$wb.Save()
